# Add two new columns (I0 -> column I, IF -> column J) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the existing headers (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for rows 2..62, column I ("I0") and column J ("IF")
$iValues = @(9,9,9,9,9,9,9,9,9,8,9,9,9,9,9,9,9,9,8,9,9,9,9,8,9,8,9,9,9,9,9,9,9,9,9,9,9,9,9,8,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,7,6,4,5,1)
$jValues = @(9,9,9,9,9,10,9,9,9,9,9,9,9,10,9,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,10,9,9,9,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,7,6,4,5,1)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
